{"js": "// Word adds a `w:lastRenderedPageBreak` marker into the \"6 Bronvermelding\"\n// heading run (content reflowed onto a new page after edits elsewhere in\n// the document), and two now-superfluous blank spacer paragraphs\n// immediately above that heading are removed.\n\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items/text,items/style\");\nawait context.sync();\n\n// Locate the \"6 Bronvermelding\" heading paragraph (the actual heading,\n// styled \"No Spacing\" / Geenafstand - not the earlier table-of-contents\n// line that references the same text).\nlet headingIndex = -1;\nfor (let i = 0; i < paras.items.length; i++) {\n  const para = paras.items[i];\n  if (para.text === \"6 Bronvermelding\" && para.style === \"No Spacing\") {\n    headingIndex = i;\n    break;\n  }\n}\nif (headingIndex === -1) {\n  throw new Error('Could not locate the \"6 Bronvermelding\" heading paragraph.');\n}\nconst heading = paras.items[headingIndex];\n\n// Remove the two blank \"No Spacing\" spacer paragraphs directly above the\n// heading.\nconst blanksToRemove = [];\nfor (let i = headingIndex - 1, found = 0; i >= 0 && found < 2; i--, found++) {\n  const candidate = paras.items[i];\n  if (candidate.text !== \"\" || candidate.style !== \"No Spacing\") break;\n  blanksToRemove.push(candidate);\n}\nblanksToRemove.forEach((p) => p.delete());\nawait context.sync();\n\n// Pull the heading paragraph's current OOXML so its paragraph-level\n// identity attributes (w14:paraId, rsids, ...) survive the rewrite, then\n// splice a <w:lastRenderedPageBreak/> in immediately before the run's\n// <w:t>, inside the same run (matching what Word itself records when it\n// re-renders/re-paginates the document).\nconst ooxml = heading.getOoxml();\nawait context.sync();\n\nconst xml = ooxml.value;\nconst pMatch = xml.match(/<w:p [^>]*w14:paraId=\"1063515B\"[^>]*>[\\s\\S]*?<\\/w:p>/);\nconst paragraphXml = pMatch ? pMatch[0] : xml.match(/<w:p\\b[^>]*>[\\s\\S]*?<\\/w:p>/)[0];\n\nconst updatedParagraphXml = paragraphXml.replace(\n  /(<w:r>\\s*<w:rPr>[\\s\\S]*?<\\/w:rPr>)(\\s*<w:t>)/,\n  \"$1<w:lastRenderedPageBreak/>$2\"\n);\n\nconst packageXml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n  \"<w:body>\" +\n  updatedParagraphXml +\n  \"</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\";\n\nconst fullRange = heading.getRange();\nfullRange.insertOoxml(packageXml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word adds a `w:lastRenderedPageBreak` marker into the \"6 Bronvermelding\"\n# heading run (content reflowed onto a new page after edits elsewhere in\n# the document), and two now-superfluous blank spacer paragraphs\n# immediately above that heading are removed.\n\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n\n# Locate the \"6 Bronvermelding\" heading paragraph (the actual heading,\n# styled \"No Spacing\" / Geenafstand - not the earlier table-of-contents\n# line that references the same text).\n$headingIndex = -1\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $p = $paras.Item($i)\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"6 Bronvermelding\" -and $p.Style.NameLocal -eq \"No Spacing\") {\n        $headingIndex = $i\n        break\n    }\n}\nif ($headingIndex -eq -1) {\n    throw \"Could not locate the '6 Bronvermelding' heading paragraph.\"\n}\n\n# Remove the two blank \"No Spacing\" spacer paragraphs directly above the\n# heading.\n$removed = 0\nwhile ($removed -lt 2) {\n    $currentParas = $d.Paragraphs\n    $candidateIndex = $headingIndex - 1\n    if ($candidateIndex -lt 1) { break }\n    $candidate = $currentParas.Item($candidateIndex)\n    $candidateText = $candidate.Range.Text.TrimEnd([char]13, [char]7)\n    if ($candidateText -ne \"\" -or $candidate.Style.NameLocal -ne \"No Spacing\") {\n        break\n    }\n    $candidate.Range.Delete()\n    $headingIndex = $headingIndex - 1\n    $removed = $removed + 1\n}\n\n# Pull the heading paragraph's current OOXML so its paragraph-level\n# identity attributes (w14:paraId, rsids, ...) survive the rewrite, then\n# splice a <w:lastRenderedPageBreak/> in immediately before the run's\n# <w:t>, inside the same run (matching what Word itself records when it\n# re-renders/re-paginates the document).\n$headingParas = $d.Paragraphs\n$heading = $headingParas.Item($headingIndex)\n$headingRange = $heading.Range\n\n$fullXml = $headingRange.WordOpenXML\n$paraMatch = [regex]::Match($fullXml, '<w:p\\b[\\s\\S]*?</w:p>')\n$paragraphXml = $paraMatch.Value\n\n$updatedParagraphXml = $paragraphXml -replace '(<w:r>\\s*<w:rPr>[\\s\\S]*?</w:rPr>)(\\s*<w:t>)', '$1<w:lastRenderedPageBreak/>$2'\n\n$packageXml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' `\n    + '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' `\n    + '<pkg:xmlData>' `\n    + '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' `\n    + '<w:body>' + $updatedParagraphXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$headingRange.InsertXML($packageXml)\n"}
